# Update countries & provincias Spain
#
# 1) Swap the "Lugo" / "Almeria" rows (their shared-string order changed in
#    the sst, which - combined with the original <v> indices staying put on
#    the cells - means the row that used to read "Lugo" now reads "Almeria"
#    and vice versa; the "Casos activos" counts (column C) travel with the
#    labels: Lugo=5, Almeria=72 become Almeria=72 (row47), Lugo=5 (row48).
# 2) Bump the "Datos actualizados" timestamp from 06:16 to 06:46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 was Lugo (Casos activos = 5); becomes Almeria (Casos activos = 72)
$ws.Range("A47").Value = "Almeria"
$ws.Range("C47").Value = 72

# Row 48 was Almeria (Casos activos = 72); becomes Lugo (Casos activos = 5)
$ws.Range("A48").Value = "Lugo"
$ws.Range("C48").Value = 5

# Update the "last refreshed" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 06:46"
